$wb = $excel.ActiveWorkbook
$ws1 = $wb.ActiveSheet
$ws1.Name = "PLANEAMIENTO"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "ESTRUCTURA_TRABAJO_US"

# --- Sheet1: PLANEAMIENTO ---
$ws1.Range("B2").Value = "ID_US"
$ws1.Range("C2").Value = "USER_STORY"
$ws1.Range("B3").Select()

# --- Sheet2: ESTRUCTURA_TRABAJO_US ---

# Column widths (approximate target widths; runtime snaps to a 1/6 character grid)
$ws2.Columns.Item(1).ColumnWidth = 4.833333333333334
$ws2.Columns.Item(2).ColumnWidth = 20.833333333333336
$ws2.Columns.Item(3).ColumnWidth = 40.0
$ws2.Columns.Item(4).ColumnWidth = 70.16666666666667
$ws2.Columns.Item(5).ColumnWidth = 62.5

# Row 6 height
$ws2.Rows.Item(6).RowHeight = 195

# helper color: White, Background 1, Darker 25% (#BFBFBF)
$grayFill = 12566463

# Values + base alignment (left/top, no wrap) for plain body cells
$bodyAddrs = @("B3","B4","B5","B8","B9","C3","C5","C6")
$bodyValues = @("Requerimiento","Analisis","Diseño","Prueba Unitaria","Implementacion","Especificar UserStory <ID>","Diseñar Lógica de Aplicación","Diseñar Tablas")
for ($i = 0; $i -lt $bodyAddrs.Length; $i++) {
    $c = $ws2.Range($bodyAddrs[$i])
    $c.Value = $bodyValues[$i]
    $c.HorizontalAlignment = -4131
    $c.VerticalAlignment = -4160
}

# Header cells without wrap (B2:C2) - bold, gray fill, left/top align
$ws2.Range("B2").Value = "DISCIPLINA"
$ws2.Range("C2").Value = "TAREAS"
$hdrNoWrap = $ws2.Range("B2:C2")
$hdrNoWrap.Font.Bold = $true
$hdrNoWrap.Interior.Color = $grayFill
$hdrNoWrap.HorizontalAlignment = -4131
$hdrNoWrap.VerticalAlignment = -4160

# Header cells with wrap (D2:E2) - bold, gray fill, left/top align, wrap text
$ws2.Range("D2").Value = "PASOS A SEGUIR"
$ws2.Range("E2").Value = "LINEAMIENTOS"
$hdrWrap = $ws2.Range("D2:E2")
$hdrWrap.Font.Bold = $true
$hdrWrap.Interior.Color = $grayFill
$hdrWrap.HorizontalAlignment = -4131
$hdrWrap.VerticalAlignment = -4160
$hdrWrap.WrapText = $true

# Body cells with wrap (D6:E7) - left/top align, wrap text
$ws2.Range("D6").Value = "1. Diseñar las tablas requeridas para soportar las necesidades de información de la User Story (/implementacion/worklist/worklist-dac / src / main / resources / db / migration"
$ws2.Range("D7").Value = "2. Se debe actualizar el Diccionario de Datos con los cambios introducidos."
$ws2.Range("E7").Value = "El archivo se encuentra ubicado en / Diseno / Diccionario_Datos.docx"

$bodyWrap = $ws2.Range("D6:E7")
$bodyWrap.HorizontalAlignment = -4131
$bodyWrap.VerticalAlignment = -4160
$bodyWrap.WrapText = $true

# Rich text cell E6 ("Nota 1: ..." with bold red run and bold ".sql" run)
$e6 = $ws2.Range("E6")
$e6.Value = "Nota 1: Se debe revisar si las tablas en la BD ya existen asi como los campos de datos y relaciones con otras tablas. `nSe debe colocar en archivo con extension .sql siguiente el formato para su nombrado siguiente.`na) Para crear o modificar la estructura de la base de datos se debera utilizar:`nV4_(consecutivo)_Create_TaskManager_Schema.sql`nb) Si se colocaran registros como ejemplos en las tablas se debera utilizar:`nV4_(consecutivo)_Fill_TaskManager_Default.sql"
$e6.HorizontalAlignment = -4131
$e6.VerticalAlignment = -4160
$e6.WrapText = $true

Write-Host "Sheet2 content done"

$ws2.Select()
$ws2.Range("E8").Select()

Write-Host "All done"
